$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.552.55"
$ws.Range("E2").Value = "  -0.79%  "

# Row 3
$ws.Range("D3").Value = "3.542.76"
$ws.Range("E3").Value = "  -2.11%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").Value = "'197.57"
$ws.Range("E5").Value = "  +0.97%  "

# Row 6
$ws.Range("D6").Value = "'586.45"
$ws.Range("E6").Value = "  -3.05%  "

# Row 7
$ws.Range("D7").Value = "'0.614"
$ws.Range("E7").Value = "  -1.91%  "

# Row 8
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").Value = "'0.208"
$ws.Range("E9").Value = "  +1.00%  "

# Row 10
$ws.Range("D10").Value = "'0.629"
$ws.Range("E10").Value = "  -3.42%  "

# Row 11
$ws.Range("D11").Value = "'52.09"
$ws.Range("E11").Value = "  -3.41%  "

# Row 12
$ws.Range("D12").Value = "'0.0000289"
$ws.Range("E12").Value = "  -4.88%  "

# Row 13
$ws.Range("D13").Value = "'9.36"
$ws.Range("E13").Value = "  -1.89%  "

# Row 14
$ws.Range("D14").Value = "'680.21"
$ws.Range("E14").Value = "  +14.84%  "

# Row 15
$ws.Range("D15").Value = "4.096.67"
$ws.Range("E15").Value = "  -2.20%  "

# Row 16
$ws.Range("D16").Value = "69.598.51"
$ws.Range("E16").Value = "  -1.00%  "

# Row 17
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "'12.47"
$ws.Range("E17").Value = "  -6.00%  "

# Row 18
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'18.63"
$ws.Range("E18").Value = "  -3.44%  "

# Row 19
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.524.19"
$ws.Range("E19").Value = "  -1.94%  "

# Row 20
$ws.Range("E20").Value = "  -0.84%  "

# Row 21
$ws.Range("D21").Value = "'0.971"
$ws.Range("E21").Value = "  -2.62%  "

# Row 22
$ws.Range("D22").Value = "'17.95"
$ws.Range("E22").Value = "  +0.30%  "

# Row 23
$ws.Range("D23").Value = "'107.84"
$ws.Range("E23").Value = "  +4.75%  "

# Row 24
$ws.Range("D24").Value = "'5.27"
$ws.Range("E24").Value = "  +2.08%  "

# Row 25
$ws.Range("D25").Value = "'4.42"
$ws.Range("E25").Value = "  -4.89%  "

# Row 26
$ws.Range("D26").Value = "'2.96"
$ws.Range("E26").Value = "  -3.69%  "

# Row 27
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'6.00"
$ws.Range("E27").Value = "  -0.94%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.36"
$ws.Range("E28").Value = "  -4.76%  "

# Row 29
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'9.73"
$ws.Range("E29").Value = "  +1.29%  "

# Row 30
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'33.48"
$ws.Range("E30").Value = "  -1.42%  "

# Row 31
$ws.Range("B31").Value = "dogwifhat"
$ws.Range("C31").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D31").Value = "'4.41"
$ws.Range("E31").Value = "  +0.35%  "

# Row 32
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'6.94"
$ws.Range("E32").Value = "  -2.66%  "

# Row 33
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'11.96"
$ws.Range("E33").Value = "  -3.06%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.112"
$ws.Range("E34").Value = "  -3.95%  "

# Row 35
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'62.34"
$ws.Range("E35").Value = "  -1.44%  "

# Row 36
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.804.77"
$ws.Range("E36").Value = "  -3.22%  "

# Row 37
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0821"
$ws.Range("E37").Value = "  -4.57%  "

# Row 38
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.12%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'3.62"
$ws.Range("E39").Value = "  +2.16%  "

# Row 40
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'505.02"
$ws.Range("E40").Value = "  -3.74%  "

# Row 41
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "'2.97"
$ws.Range("E41").Value = "  -6.64%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.136"
$ws.Range("E42").Value = "  +1.77%  "

# Row 43
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "'0.374"
$ws.Range("E43").Value = "  -4.95%  "

# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'34.98"
$ws.Range("E44").Value = "  -5.86%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0461"
$ws.Range("E45").Value = "  +0.68%  "

# Row 46
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "'2.97"
$ws.Range("E46").Value = "  +3.90%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.37"
$ws.Range("E47").Value = "  +1.42%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.138"
$ws.Range("E48").Value = "  -2.34%  "

# Row 49
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  -0.46%  "

# Row 50
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'8.39"
$ws.Range("E50").Value = "  -2.69%  "

# Row 51
$ws.Range("B51").Value = "Jupiter"
$ws.Range("C51").Value = "https://coinranking.com/coin/qMgTxtv34+jupiter-jup"
$ws.Range("D51").Value = "'1.82"
$ws.Range("E51").Value = "  +22.15%  "
